# Investimentos.xlsx edit script
# 1) Replace row 2 data (mglu3 -> petr4) and append rows 3-5 (vale3, cmig3, azul4)
# 2) Recolor the 3 shared conditional-format fills (green / red / white)
# 3) Rebuild conditional formatting across F, Q, N, R, O, K, L, S, M, P, I (rows 2:5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Data
# ---------------------------------------------------------------------------
$rows = @{
    2 = @("petr4", @(38.62, 31.25, 72.09066666666666, 79.50530642667822, 23.67, 8.99, 16.8, 23.27809425168307, 80.91662351113413, 4.3, 1.24, 28.77, 19.56, 100, 2.89, 0.5600000000000001, 0.88, 1.08))
    3 = @("vale3", @(62.7, 41.24, 76.91199999999999, 88.96621268773893, 18.78, 8.529999999999999, 11.04, 13.60446570972887, 65.77352472089314, 7.42, 1.53, 20.68, 18, 100, 4.35, 0.29, 0.6899999999999999, 1.12))
    4 = @("cmig3", @(12.41, 8.880000000000001, 14.87821111111111, 19.63705680594727, 14.81, 1.93, 10.79, 15.55197421434327, 71.55519742143433, 6.53, 1.42, 21.72, 16.13, 80, 4.63, 0.28, 0.86, 1.15))
    5 = @("azul4", @(7.12, -17.69, 0, 0, 10.81, 1.6, 0, 22.47191011235955, -248.4550561797753, 4.55, -0.41, -9.06, 43.38, 100, 0, 0, 0, 0.34))
}

foreach ($r in 2..5) {
    $entry = $rows[$r]
    $name = $entry[0]
    $vals = $entry[1]
    $ws.Cells.Item($r, 1).Value = $name
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------------
# 2) Recolor the 3 shared dxf fills used by conditional formatting
#    (green = 92D050, red = FF0000, white = FFFFFF, as VB BGR-packed longs)
# ---------------------------------------------------------------------------
$green = 5296274    # RGB(146,208,80)  0x92D050
$red   = 255         # RGB(255,0,0)     0xFF0000
$white = 16777215    # RGB(255,255,255) 0xFFFFFF

$f2conds = $ws.Range("F2").FormatConditions
$fGreen = $f2conds.Item(1)
$fRed   = $f2conds.Item(2)
$fWhite = $f2conds.Item(3)

$fGreen.Interior.Color = $green
$fRed.Interior.Color   = $red
$fWhite.Interior.Color = $white

# ---------------------------------------------------------------------------
# 3) Rebuild the F2 rule (>=10 / between / equal-empty) over F2:F5
# ---------------------------------------------------------------------------
$fGreen.Operator = 7
$fGreen.Formula1 = "10"
$fGreen.Priority = 1

$fRed.Operator = 1
$fRed.Formula1 = "0.001"
$fRed.Formula2 = "9.999"
$fRed.Priority = 2

$fWhite.Operator = 3
$fWhite.Formula1 = ""
$fWhite.Priority = 3

$fGreen.ModifyAppliesToRange($ws.Range("F2:F5"))

# ---------------------------------------------------------------------------
# 4) Build the remaining 10 conditional-formatting ranges.
#    Operator codes (XlFormatConditionOperator): xlBetween=1, xlEqual=3,
#    xlGreaterEqual=7
# ---------------------------------------------------------------------------
function Add-Rule($range, $op, $f1, $f2, $color, $priority) {
    if ($f2 -eq $null) {
        $c = $range.FormatConditions.Add(1, $op, $f1)
    } else {
        $c = $range.FormatConditions.Add(1, $op, $f1, $f2)
    }
    $c.Interior.Color = $color
    $c.Priority = $priority
    return $c
}

# Q2:Q5  -> between(0.001,0.999)=red, >=1=green, equal(empty)=white
$rngQ = $ws.Range("Q2:Q5")
Add-Rule $rngQ 1 "0.001" "0.999" $red   4  | Out-Null
Add-Rule $rngQ 7 "1"     $null   $green 5  | Out-Null
Add-Rule $rngQ 3 ""      $null   $white 6  | Out-Null

# N2:N5  -> >=15=green, between(0.001,14.999)=red, equal(empty)=white
$rngN = $ws.Range("N2:N5")
Add-Rule $rngN 7 "15"    $null     $green 7  | Out-Null
Add-Rule $rngN 1 "0.001" "14.999"  $red   8  | Out-Null
Add-Rule $rngN 3 ""      $null     $white 9  | Out-Null

# R2:R5  -> between(0.001,2.999)=red, >=3=green, equal(empty)=white
$rngR = $ws.Range("R2:R5")
Add-Rule $rngR 1 "0.001" "2.999" $red   10 | Out-Null
Add-Rule $rngR 7 "3"     $null   $green 11 | Out-Null
Add-Rule $rngR 3 ""      $null   $white 12 | Out-Null

# O2:O5  -> equal(100)=green, equal(empty)=white  (no red rule)
$rngO = $ws.Range("O2:O5")
Add-Rule $rngO 3 "100" $null   $green 13 | Out-Null
Add-Rule $rngO 3 ""    $null   $white 14 | Out-Null

# K2:K5  -> between(0.001,9.999)=red, >=10=green, equal(empty)=white
$rngK = $ws.Range("K2:K5")
Add-Rule $rngK 1 "0.001" "9.999" $red   15 | Out-Null
Add-Rule $rngK 7 "10"    $null   $green 16 | Out-Null
Add-Rule $rngK 3 ""      $null   $white 17 | Out-Null

# L2:L5  -> between(0.001,1.499)=red, >=1.5=green, equal(empty)=white
$rngL = $ws.Range("L2:L5")
Add-Rule $rngL 1 "0.001" "1.499" $red   18 | Out-Null
Add-Rule $rngL 7 "1.5"   $null   $green 19 | Out-Null
Add-Rule $rngL 3 ""      $null   $white 20 | Out-Null

# S2:S5  -> >=1=green, between(0.001,0.999)=red, equal(empty)=white
$rngS = $ws.Range("S2:S5")
Add-Rule $rngS 7 "1"     $null     $green 21 | Out-Null
Add-Rule $rngS 1 "0.001" "0.999"   $red   22 | Out-Null
Add-Rule $rngS 3 ""      $null     $white 23 | Out-Null

# M2:M5  -> >=16=green, between(0.001,15.999)=red, equal(empty)=white
$rngM = $ws.Range("M2:M5")
Add-Rule $rngM 7 "16"    $null      $green 24 | Out-Null
Add-Rule $rngM 1 "0.001" "15.999"   $red   25 | Out-Null
Add-Rule $rngM 3 ""      $null      $white 26 | Out-Null

# P2:P5  -> between(0.001,4.999)=red, >=5=green, equal(empty)=white
$rngP = $ws.Range("P2:P5")
Add-Rule $rngP 1 "0.001" "4.999" $red   27 | Out-Null
Add-Rule $rngP 7 "5"     $null   $green 28 | Out-Null
Add-Rule $rngP 3 ""      $null   $white 29 | Out-Null

# I2:I5  -> >=20=green, between(0.001,19.999)=red, equal(empty)=white
$rngI = $ws.Range("I2:I5")
Add-Rule $rngI 7 "20"    $null      $green 30 | Out-Null
Add-Rule $rngI 1 "0.001" "19.999"   $red   31 | Out-Null
Add-Rule $rngI 3 ""      $null      $white 32 | Out-Null

Write-Output "done"
